$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (order chosen so shared-string table indices
# line up with the saved file: 0=Left, 1=Others, 2=Right)
$ws.Range("A1").Value = "Left Marking Text (L)"
$ws.Range("C1").Value = "Others (O)"
$ws.Range("B1").Value = "Right Marking Text  (R)"

# Remove the sample data rows (2 and 3)
$ws.Range("A2:C3").ClearContents()

# Update selection to match the saved file state
$ws.Range("B3").Select()
